# "fixed bug for swap scene"
# The Scene config sheet has two boolean-like columns:
#   J = "Share", K = "CanClone"
# For rows 11-35 the K ("CanClone") flag was stuck at 0 while it should be 1.
# Row 12 additionally had J and K swapped (J=1/K=0 instead of J=0/K=1).
# This script fixes those cell values and restores the worksheet selection
# to the single cell K13 (matching the state the workbook was left in after
# the fix was applied and verified).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 11; $r -le 35; $r++) {
    if ($r -eq 12) {
        $ws.Range("J$r").Value = 0
    }
    $ws.Range("K$r").Value = 1
}

$ws.Range("K13").Select()
